$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new rows (15-17) with the delete-transaction API entries.
# Column B (API keys) is populated first, then column C (descriptions),
# matching the order the new entries were authored.
$ws.Range("B15").Value = "transaction.delete.master.setBusinessDocumentType"
$ws.Range("B17").Value = "transaction.delete.master.setCountry"
$ws.Range("B16").Value = "transaction.delete.master.setCitizenIdentity"

$ws.Range("C15").Value = "Menghapus Data Jenis Dokumen Bisnis"
$ws.Range("C16").Value = "Menghapus Data Identitas Penduduk"
$ws.Range("C17").Value = "Menghapus Data Negara"

# Move the active selection to the bottom-right pane's new last-entered cell
$ws.Range("C18").Select()
